# Set the "Diferencia Stock" (column L) values to 0 for all rows that
# currently hold a non-zero adjustment, and reset the
# "Total_Ajuste_Stock:" summary cell (C107) to 0 to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Semana_8")

$rows = @(5, 8, 12, 15, 17, 19, 22, 23, 26, 31, 33, 34, 35, 37, 38, 39, 40, 41, 42, 43, 46, 49, 51, 53, 55, 56, 58, 60, 61, 62, 64, 65, 66, 69, 70, 71, 77, 85, 90, 93)

foreach ($r in $rows) {
    $ws.Range("L$r").Value = 0
}

$ws.Range("C107").Value = 0
